# Applies the "added tasks Read Feature in projects page (CRUD)" commit:
#  - Marks "Task 15: Tasks Factory Seeder" and "Task 16: Projects Read (CRUD)" as
#    Complete and stamps them with the 02/28/2024 date.
#  - Marks "Task 19: Projects Delete (CRUD)" as Pending.
#  - Fills in the previously-blank "Task 20:".."Task 23:" rows with their full
#    descriptions (Tasks Read/Create/Update/Delete - CRUD) and statuses.
#  - Leaves the selection on D23, matching the author's last-edited cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextDate {
    param($Cell, $Text)
    # Writing a date-shaped string straight into Value2 gets auto-coerced into a
    # serial date number by Excel. Routing it through a text formula first and
    # then collapsing that formula down to a static value (paste-values) keeps
    # it as a genuine shared-string/text cell instead, matching the workbook.
    $Cell.Formula = '="' + $Text + '"'
    $Cell.Copy()
    $Cell.PasteSpecial(-4163)
}

# --- Row 15: Task 15: Tasks Factory Seeder ---
$ws.Range("B15").Value2 = "Complete"
Set-TextDate $ws.Range("D15") "02/28/2024"

# --- Row 16: Task 16: Projects Read (CRUD) ---
$ws.Range("B16").Value2 = "Complete"
Set-TextDate $ws.Range("D16") "02/28/2024"

# --- Row 19: Task 19: Projects Delete (CRUD) ---
$ws.Range("B19").Value2 = "Pending"

# --- Row 20: Task 20: Tasks Read (CRUD) ---
$ws.Range("A20").Value2 = "Task 20: Tasks Read (CRUD)"
$ws.Range("B20").Value2 = "Complete"
Set-TextDate $ws.Range("D20") "02/28/2024"

# --- Row 21: Task 21: Tasks Create (CRUD) ---
$ws.Range("A21").Value2 = "Task 21: Tasks Create (CRUD)"
$ws.Range("B21").Value2 = "Pending"

# --- Row 22: Task 22: Tasks Update (CRUD) ---
$ws.Range("A22").Value2 = "Task 22: Tasks Update (CRUD)"
$ws.Range("B22").Value2 = "Pending"

# --- Row 23: Task 23: Tasks Delete (CRUD) ---
$ws.Range("A23").Value2 = "Task 23: Tasks Delete (CRUD)"
$ws.Range("B23").Value2 = "Pending"
Set-TextDate $ws.Range("D23") " "

# Now re-apply the correct status colour/formatting to every Status cell by
# copying the format from an existing cell that already carries the right
# style (Complete -> green, Pending -> orange), without touching the values.
$ws.Range("B2").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("B20").PasteSpecial(-4122)

$ws.Range("B17").Copy()
$ws.Range("B19").PasteSpecial(-4122)
$ws.Range("B21").PasteSpecial(-4122)
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("B23").PasteSpecial(-4122)

# Clear clipboard marching-ants state and restore the last-used selection.
$excel.CutCopyMode = $false
$ws.Range("D23").Select() | Out-Null
